## vault backup edit: duplicate "Mdm权限管理" as a "-bak" copy placed right
## before it, then trim the original sheet's data down to just its header
## rows (the commit is a routine Obsidian vault sync of the authors' working
## copy of the permission list; the functional change inside the workbook is
## that the full "Mdm权限管理" sheet was backed up and then had its data body
## (rows 5-74) removed).

$wb = $excel.ActiveWorkbook

# 1) Duplicate "Mdm权限管理" -> Excel inserts the copy immediately before the
#    source sheet when Before:=source, After:=nothing; use Copy(Before) which
#    is the "Move or Copy... Create a copy" behaviour.
$src = $wb.Worksheets.Item("Mdm权限管理")
$src.Copy($src)

# 2) The freshly inserted copy now sits directly before the (re-resolved)
#    original sheet. Rename it to the "-bak" name used by the author.
$orig = $wb.Worksheets.Item("Mdm权限管理")
$bak = $wb.Worksheets.Item($orig.Index - 1)
$bak.Name = "Mdm权限管理-bak"

# 3) Trim the live "Mdm权限管理" sheet back down to just the header block
#    (rows 1-4); the rest of the permission rows now only live in the backup
#    copy.
$orig.Rows("5:74").Delete()

# 4) Restore the selection / active-sheet cosmetics that came along with the
#    edit: the backup keeps the old selection, the live sheet is the newly
#    active tab with a fresh selection.
$bak.Range("H30").Select()
$orig.Activate()
$orig.Range("F18").Select()
